$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 757
$ws1.Range("F3").Value = 14
$ws1.Range("F5").Value = 25
$ws1.Range("F7").Value = 3438
$ws1.Range("F8").Value = 71
$ws1.Range("F9").Value = 4110
$ws1.Range("F11").Value = 1025

# Sheet "全部类型" - same underlying records, one extra row offset
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 757
$ws4.Range("F3").Value = 14
$ws4.Range("F5").Value = 25
$ws4.Range("F8").Value = 3438
$ws4.Range("F9").Value = 71
$ws4.Range("F10").Value = 4110
$ws4.Range("F12").Value = 1025
